# Add new columns I (I0) and J (IF) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold, border, centered) from H1 onto the new
# header cells I1:J1 before assigning their text, so the style survives.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new data columns for rows 2-14.
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 7

$ws.Range("I3").Value = 7
$ws.Range("J3").Value = 7

$ws.Range("I4").Value = 5
$ws.Range("J4").Value = 5

$ws.Range("I5").Value = 8
$ws.Range("J5").Value = 9

$ws.Range("I6").Value = 7
$ws.Range("J6").Value = 7

$ws.Range("I7").Value = 8
$ws.Range("J7").Value = 8

$ws.Range("I8").Value = 6
$ws.Range("J8").Value = 7

$ws.Range("I9").Value = 9
$ws.Range("J9").Value = 9

$ws.Range("I10").Value = 5
$ws.Range("J10").Value = 5

$ws.Range("I11").Value = 7
$ws.Range("J11").Value = 7

$ws.Range("I12").Value = 4
$ws.Range("J12").Value = 5

$ws.Range("I13").Value = 5
$ws.Range("J13").Value = 5

$ws.Range("I14").Value = 4
$ws.Range("J14").Value = 4
